$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "feed" (sheet1)
# ---------------------------------------------------------------------------
$feed = $wb.Worksheets.Item("feed")

# New helper cells in rows 3/4/6
$feed.Range("K3").Formula = '=$A$6*2/($B$6/$C$6*A9)'

$feed.Range("J4").Formula = '=A6*2*C6/A9/B6*B1'
$feed.Range("J4").ClearFormats()
$feed.Range("K4").Formula = '=ROUNDDOWN($A$6*2/($B$6/$C$6*A9)*$B$1,0)'

$feed.Range("J6").Formula = '=IFERROR(CONCATENATE("0x",DEC2HEX(ROUNDDOWN($A$6*2/($B$6/$C$6*A9)*$B$1,0),9)),"0x000000000")'
$feed.Range("J6").HorizontalAlignment = -4108

# Remove the old B9 helper cell (its formula is now computed via J6 instead)
$feed.Range("B9").Clear()

# I9 now reads the hex string from J6 rather than from the deleted B9
$feed.Range("I9").Formula = '=$A$2&J6&$C$2&C9&$D$2&SUBSTITUTE(TEXT(A9,"0,00"),",",".")&$G$2&D9&$E$2&E9&$D$2&F9&$G$2&G9&$E$2&H9&$B$2'

# Select feed's B15 and make it the active sheet/tab
$feed.Activate()
$feed.Range("B15").Select()
